$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New package entry appended as row 10, mirroring the existing rows
# (same package "DTDemo" / iflow1 / 1.0.0 / IFlow / 2026-02-04 combo
# already used elsewhere on the sheet).
$ws.Range("A10").Value = "DTDemo"
$ws.Range("B10").Value = "DTDemo"
$ws.Range("C10").Value = "iflow1"
$ws.Range("D10").Value = "1.0.0"
$ws.Range("E10").Value = "IFlow"

# F10 holds the text "2026-02-04" (not a real date, same as F6:F9).
# Assigning it via .Value would let Excel auto-convert the ISO-looking
# string into a date serial, which is not what the source data has.
# Copy an existing text cell with the identical value instead, so the
# text type/formatting is preserved exactly like its neighbours.
$ws.Range("F9").Copy()
$ws.Range("F10").PasteSpecial()
$excel.CutCopyMode = $false
